$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 7079
$ws1.Range("F7").Value  = 919
$ws1.Range("F8").Value  = 272
$ws1.Range("F9").Value  = 761
$ws1.Range("F10").Value = 539
$ws1.Range("F11").Value = 12
$ws1.Range("C13").Value = "杭州·恋爱告急动漫游戏展"
$ws1.Range("F13").Value = 55
$ws1.Range("G13").Value = 52
$ws1.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202401/de0lsSbz1706180254353.jpeg"
$ws1.Range("F16").Value = 2831
$ws1.Range("F17").Value = 140
$ws1.Range("F18").Value = 25
$ws1.Range("F19").Value = 193
$ws1.Range("F21").Value = 37
$ws1.Range("F22").Value = 421
$ws1.Range("F23").Value = 125
$ws1.Range("F25").Value = 105
$ws1.Range("F26").Value = 154
$ws1.Range("F27").Value = 113
$ws1.Range("F32").Value = 246
$ws1.Range("F33").Value = 361

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 1

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 191

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 191
$ws4.Range("F9").Value  = 7079
$ws4.Range("F11").Value = 919
$ws4.Range("F12").Value = 272
$ws4.Range("F13").Value = 761
$ws4.Range("F14").Value = 539
$ws4.Range("F15").Value = 12
$ws4.Range("C17").Value = "杭州·恋爱告急动漫游戏展"
$ws4.Range("F17").Value = 55
$ws4.Range("G17").Value = 52
$ws4.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202401/de0lsSbz1706180254353.jpeg"
$ws4.Range("F21").Value = 2831
$ws4.Range("F22").Value = 140
$ws4.Range("F23").Value = 25
$ws4.Range("F25").Value = 193
$ws4.Range("F27").Value = 1
$ws4.Range("F28").Value = 37
$ws4.Range("F29").Value = 421
$ws4.Range("F30").Value = 125
$ws4.Range("F32").Value = 105
$ws4.Range("F33").Value = 154
$ws4.Range("F34").Value = 113
$ws4.Range("F39").Value = 246
$ws4.Range("F40").Value = 361
